$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 (Marking row): Right -> 4, Wrong -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total row): Right -> 92, Wrong -> -2, Max text -> "90 / 112"
$ws.Range("B12").Value = 92
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "90 / 112"
